# Applies the cryptos list update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data: row number, Coin (B), Link (C), Price (D), Volume(1h) (E)
$data = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "70.950.81", "  +3.36%  "),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "3.562.73", "  +2.26%  "),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.00", "  -0.09%  "),
    @(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "583.00", "  +2.40%  "),
    @(6, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "186.72", "  +2.31%  "),
    @(7, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.631", "  +3.04%  "),
    @(8, "LidoStakedEther", "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth", "3.550.87", "  +2.09%  "),
    @(9, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.00", "  -0.03%  "),
    @(10, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.219", "  +19.60%  "),
    @(11, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.652", "  +2.48%  "),
    @(12, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "54.54", "  +1.95%  "),
    @(13, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0000316", "  +5.77%  "),
    @(14, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "9.48", "  +0.86%  "),
    @(15, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "4.125.38", "  +1.90%  "),
    @(16, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "70.962.40", "  +3.36%  "),
    @(17, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "19.22", "  +0.18%  "),
    @(18, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "3.575.28", "  +2.62%  "),
    @(19, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "12.44", "  +0.50%  "),
    @(20, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "570.76", "  +6.04%  "),
    @(21, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.120", "  +0.74%  "),
    @(22, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "1.00", "  -1.05%  "),
    @(23, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "17.66", "  -10.05%  "),
    @(24, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "4.55", "  +3.82%  "),
    @(25, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "4.94", "  -0.81%  "),
    @(26, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "95.38", "  +1.59%  "),
    @(27, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "11.29", "  +4.34%  "),
    @(28, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "2.95", "  +1.80%  "),
    @(29, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "9.16", "  +1.80%  "),
    @(30, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "32.55", "  +3.97%  "),
    @(31, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "7.23", "  +0.79%  "),
    @(32, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "12.28", "  -1.99%  "),
    @(33, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.117", "  +3.37%  "),
    @(34, "Fetch.AI", "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet", "3.44", "  +14.07%  "),
    @(35, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "63.11", "  -1.65%  "),
    @(36, "Bittensor", "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao", "548.44", "  -4.13%  "),
    @(37, "TheGraph", "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt", "0.414", "  +4.55%  "),
    @(38, "InjectiveProtocol", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj", "37.91", "  +0.37%  "),
    @(39, "dogwifhat", "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif", "3.39", "  +10.31%  "),
    @(40, "PEPE", "https://coinranking.com/coin/03WI8NQPF+pepe-pepe", "0.0₃0803", "  +5.77%  "),
    @(41, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "0.999", "  -0.03%  "),
    @(42, "Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "3.564.11", "  +11.54%  "),
    @(43, "Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "3.47", "  +4.74%  "),
    @(44, "Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.137", "  +3.18%  "),
    @(45, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.0461", "  +5.75%  "),
    @(46, "ApeXProtocol", "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex", "3.49", "  +0.93%  "),
    @(47, "ThetaToken", "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta", "2.91", "  -1.27%  "),
    @(48, "THORChain", "https://coinranking.com/coin/ybmU-kKU+thorchain-rune", "9.32", "  +2.76%  "),
    @(49, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.138", "  +3.25%  "),
    @(50, "OceanProtocol", "https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean", "1.52", "  +16.67%  "),
    @(51, "FirstDigitalUSD", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", "0.999", "  +0.00%  ")
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]

    # Price column: values like "1.00" / "0.219" must stay literal text
    # (matching the source inlineStr cells) instead of being parsed into
    # numbers by Excel's normal value-assignment coercion.
    $priceCell = $ws.Cells.Item($r, 4)
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $item[3]
    $priceCell.Style = "Normal"

    $ws.Cells.Item($r, 5).Value = $item[4]
}
